# Apply updated cryptocurrency price/volume data to the worksheet.
# Column D values are forced to Text format (matching the source inlineStr
# cells) so numeric-looking strings like "218.82" are not auto-converted
# to numbers by Excel; the style is reset back to Normal afterwards so no
# extraneous cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.785.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.645.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.69"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.764.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.16%  "

$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "211.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("E27").Value = "  -1.72%  "

$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.10%  "

$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.277.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  -0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.528"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.808"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.59%  "

$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  -2.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.779.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("E47").Value = "  -1.29%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.77%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.10%  "
